$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.415.45"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "1.935.29"
$ws.Range("E3").Value = "  -2.03%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.62%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.91%  "

$ws.Range("E9").Value = "  -3.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0841"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").Value = "2.218.38"
$ws.Range("E12").Value = "  -1.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.802"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.59%  "

$ws.Range("D17").Value = "1.929.25"
$ws.Range("E17").Value = "  -2.35%  "

$ws.Range("D18").Value = "36.360.00"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.35%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "226.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.89%  "

$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.31%  "

$ws.Range("E25").Value = "  -2.37%  "

$ws.Range("E26").Value = "  -7.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.19%  "

$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.65%  "

$ws.Range("E32").Value = "  -6.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.44%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.05%  "

$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("E38").Value = "  -3.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0971"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("E42").Value = "  -2.27%  "

$ws.Range("E43").Value = "  -5.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.54%  "

$ws.Range("D45").Value = "1.329.29"
$ws.Range("E45").Value = "  -2.41%  "

$ws.Range("E46").Value = "  -6.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.26%  "

$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").Value = "2.110.58"
$ws.Range("E50").Value = "  -1.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.70%  "
